$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# LESSON 9 (rows 99-108): fill in beginning (B), word (D), translat (F),
# ending (H) columns — completing the exercise rows that were previously
# placeholders / partially filled.
# ---------------------------------------------------------------------------

$ws.Range("B99").Value  = "She"
$ws.Range("D99").Value  = "читає"
$ws.Range("F99").Value  = "is not reading"
$ws.Range("H99").Value  = "a new book"

$ws.Range("B100").Value = "They"
$ws.Range("D100").Value = "не тримали"
$ws.Range("F100").Value = "were not keeping"
$ws.Range("H100").Value = "the door open"

$ws.Range("B101").Value = "We"
$ws.Range("D101").Value = "приєднаємося"
$ws.Range("F101").Value = "will be joining"
$ws.Range("H101").Value = "this project"

$ws.Range("B102").Value = "He"
$ws.Range("D102").Value = "починає"
$ws.Range("F102").Value = "is starting"
$ws.Range("H102").Value = "a new job"

$ws.Range("B103").Value = "I"
$ws.Range("D103").Value = "вивчав"
$ws.Range("F103").Value = "was learning"
$ws.Range("H103").Value = "a lot last year"

$ws.Range("B104").Value = "The teacher"
$ws.Range("D104").Value = "допомагатиме"
$ws.Range("F104").Value = "will be helping"
$ws.Range("H104").Value = "every student"

$ws.Range("B105").Value = "They"
$ws.Range("D105").Value = "працюють"
$ws.Range("F105").Value = "are working"
$ws.Range("H105").Value = "in the field now"

$ws.Range("B106").Value = "We"
$ws.Range("D106").Value = "будемо будувати"
$ws.Range("F106").Value = "will be building"
$ws.Range("H106").Value = "a large house"

$ws.Range("B107").Value = "She"
$ws.Range("D107").Value = "кликала"
$ws.Range("F107").Value = "was calling"
$ws.Range("H107").Value = "us"

$ws.Range("B108").Value = "He"
$ws.Range("D108").Value = "зробив"
$ws.Range("F108").Value = "was making"
$ws.Range("H108").Value = "the files available"

# ---------------------------------------------------------------------------
# LESSON 10 (rows 111-120): fill in beginning (B), translat (F) and ending
# (H) columns. The word (D) column is left as the original placeholder text
# for this lesson.
# ---------------------------------------------------------------------------

$ws.Range("B111").Value = "I"
$ws.Range("F111").Value = "have lived"
$ws.Range("H111").Value = "here for 10 years"

$ws.Range("B112").Value = "We both"
$ws.Range("F112").Value = "have worked"
$ws.Range("H112").Value = "here since Summer"

$ws.Range("B113").Value = "They"
$ws.Range("F113").Value = "have given"
$ws.Range("H113").Value = "us the keys"

$ws.Range("B114").Value = "Scientists"
$ws.Range("F114").Value = "have studied"
$ws.Range("H114").Value = "human behaviour for centuries"

$ws.Range("B115").Value = "Both teams"
$ws.Range("F115").Value = "have played"
$ws.Range("H115").Value = "well"

$ws.Range("B116").Value = "She"
$ws.Range("F116").Value = "has searched"
$ws.Range("H116").Value = "for it since last week"

$ws.Range("B117").Value = "We"
$ws.Range("F117").Value = "have lived"
$ws.Range("H117").Value = "without the car for a month"

$ws.Range("B118").Value = "She"
$ws.Range("F118").Value = "has found"
$ws.Range("H118").Value = "something interesting"

$ws.Range("B119").Value = "They"
$ws.Range("F119").Value = "have come"
$ws.Range("H119").Value = "to visit us"

$ws.Range("B120").Value = "We"
$ws.Range("F120").Value = "have met"
$ws.Range("H120").Value = "the local guide"

# ---------------------------------------------------------------------------
# Update the sheet's active selection to match where the author left off
# (scrolled down near the just-completed Lesson 10 block, cell J104 active).
# ---------------------------------------------------------------------------

$ws.Range("J104").Select() | Out-Null
